{"js": "// Office.js (Word JavaScript API) edit script\n// Replaces the text of the first 7 paragraphs and appends 5 new paragraphs\n// at the end, matching the diff in the commit.\n\n// Paragraph 1's two runs (separated by a manual line break <w:br/>) are\n// represented in Office.js as a single `text` string joined with the\n// vertical-tab character (\\u000b) that stands in for the line break.\nconst existingParagraphTexts = [\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 - 14.03.25\\u000bA Survey on Kolmogorov-Arnold Network\",\n  \"\u05de\u05d1\u05d5\u05d0:\",\n  \" \u05d6\u05d5\u05db\u05e8\u05d9\u05dd \u05d0\u05ea KANs? \u05e9\u05d6\u05d4 \u05e7\u05d9\u05e6\u05d5\u05e8 \u05e9\u05dc Kolmogorov-Arnold Networks \u05e9\u05e2\u05e9\u05d4 \u05d4\u05e8\u05d1\u05d4 \u05e8\u05e2\u05e9 \u05d1\u05d6\u05de\u05e0\u05d5 \u05d0\u05da \u05d4\u05d1\u05d0\u05d6 \u05d4\u05dc\u05da \u05d5\u05d3\u05e2\u05da \u05e2\u05dd \u05d4\u05d6\u05de\u05df. \u05de\u05ea\u05d1\u05e8\u05e8 \u05e9\u05d9\u05e6\u05d0\u05d5 \u05dc\u05d0 \u05de\u05e2\u05d8 \u05de\u05d7\u05e7\u05e8\u05d9\u05dd \u05d1\u05e0\u05d5\u05e9\u05d0 \u05d4\u05de\u05e8\u05ea\u05e7 \u05d4\u05d6\u05d4. \u05d4\u05de\u05d0\u05de\u05e8 \u05d3\u05df \u05d1\u05d4\u05e8\u05d7\u05d1\u05d5\u05ea \u05d5\u05e9\u05d9\u05e0\u05d5\u05d9\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd \u05dc\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05ea \u05d4-KAN \u05d4\u05d1\u05e1\u05d9\u05e1\u05d9\u05ea. \u05d0\u05dc\u05d4 \u05db\u05d5\u05dc\u05dc\u05d9\u05dd \u05d4\u05ea\u05d0\u05de\u05d5\u05ea \u05dc\u05e0\u05d9\u05ea\u05d5\u05d7 \u05e1\u05d3\u05e8\u05d5\u05ea \u05e2\u05ea\u05d9\u05d5\u05ea, \u05dc\u05e2\u05d9\u05d1\u05d5\u05d3 \u05d3\u05d0\u05d8\u05d4 \u05d2\u05e8\u05e4\u05d9 \u05d5\u05dc\u05e4\u05ea\u05e8\u05d5\u05df \u05de\u05e9\u05d5\u05d5\u05d0\u05d5\u05ea \u05d3\u05d9\u05e4\u05e8\u05e0\u05e6\u05d9\u05d0\u05dc\u05d9\u05d5\u05ea. \u05e9\u05d9\u05e0\u05d5\u05d9\u05d9\u05dd \u05d0\u05dc\u05d4 \u05db\u05d5\u05dc\u05dc\u05d9\u05dd \u05dc\u05e8\u05d5\u05d1 \u05e9\u05d9\u05dc\u05d5\u05d1 \u05e9\u05dc \u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05de\u05d9\u05d5\u05d7\u05d3\u05d9\u05dd \u05d0\u05d5 \u05d0\u05d9\u05dc\u05d5\u05e6\u05d9\u05dd \u05d1\u05ea\u05d5\u05da \u05d4-KAN \u05d1\u05de\u05d8\u05e8\u05d4 \u05dc\u05d4\u05ea\u05de\u05d5\u05d3\u05d3 \u05d8\u05d5\u05d1 \u05d9\u05d5\u05ea\u05e8 \u05e2\u05dd \u05d4\u05d3\u05e8\u05d9\u05e9\u05d5\u05ea \u05d4\u05e1\u05e4\u05e6\u05d9\u05e4\u05d9\u05d5\u05ea \u05e9\u05dc \u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd \u05d0\u05dc\u05d4.\",\n  \"\u05e8\u05e9\u05ea\u05d5\u05ea \u05e7\u05d5\u05dc\u05de\u05d5\u05d2\u05d5\u05e8\u05d5\u05d1-\u05d0\u05e8\u05e0\u05d5\u05dc\u05d3 \u05de\u05d9\u05d9\u05e6\u05d2\u05d5\u05ea \u05e9\u05d9\u05e0\u05d5\u05d9 \u05e4\u05e8\u05d3\u05d9\u05d2\u05de\u05d4 \u05d1\u05ea\u05db\u05e0\u05d5\u05df \u05e8\u05e9\u05ea\u05d5\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd, \u05d4\u05de\u05d1\u05d5\u05e1\u05e1\u05d5\u05ea \u05e2\u05dc \u05de\u05e2\u05d1\u05e8 \u05de\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d4 \u05e7\u05d1\u05d5\u05e2\u05d5\u05ea \u05dc\u05e7\u05e8\u05d0\u05ea \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d4\u05e0\u05d9\u05ea\u05e0\u05d5\u05ea \u05dc\u05dc\u05de\u05d9\u05d3\u05d4 \u05d4\u05e0\u05e7\u05e8\u05d0\u05d5\u05ea b-splines. \u05d4\u05d3\u05d1\u05e8 \u05e9\u05d0\u05d1 \u05d4\u05e9\u05e8\u05d0\u05d4 \u05de\u05de\u05e9\u05e4\u05d8 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05e7\u05d5\u05dc\u05de\u05d5\u05d2\u05d5\u05e8\u05d5\u05d1-\u05d0\u05e8\u05e0\u05d5\u05dc\u05d3, \u05d4\u05d8\u05d5\u05e2\u05df \u05e9\u05db\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 \u05e8\u05e6\u05d9\u05e4\u05d4 \u05e9\u05dc \u05de\u05e9\u05ea\u05e0\u05d9\u05dd \u05de\u05e8\u05d5\u05d1\u05d9\u05dd \u05e0\u05d9\u05ea\u05e0\u05ea \u05dc\u05d9\u05d9\u05e6\u05d5\u05d2 \u05db\u05d4\u05e8\u05db\u05d1\u05d4 \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05e9\u05dc \u05de\u05e9\u05ea\u05e0\u05d4 \u05d0\u05d7\u05d3. \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05e9\u05d9\u05de\u05d5\u05e9 \u05d1\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d4\u05de\u05d9\u05d5\u05e6\u05d2\u05d5\u05ea \u05e2\u05dc \u05d9\u05d3\u05d9 \u05e1\u05e4\u05dc\u05d9\u05d9\u05e0\u05d9\u05dd(\u05e9\u05d9\u05dc\u05d5\u05d1 \u05e9\u05dc \u05e4\u05d5\u05dc\u05d9\u05e0\u05d5\u05de\u05d9\u05dd \u05d1\u05d0\u05d9\u05e0\u05d8\u05e8\u05d5\u05d5\u05dc \u05e1\u05e4\u05d5\u05d9), KANs \u05de\u05e6\u05d9\u05e2\u05d5\u05ea \u05d2\u05de\u05d9\u05e9\u05d5\u05ea \u05de\u05e9\u05d5\u05e4\u05e8\u05ea \u05d5\u05e4\u05d5\u05d8\u05e0\u05e6\u05d9\u05d0\u05dc \u05dc\u05d3\u05d9\u05d5\u05e7 \u05d2\u05d1\u05d5\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05d1\u05e7\u05d9\u05e8\u05d5\u05d1 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea. \u05d3\u05d1\u05e8 \u05de\u05d5\u05d1\u05d9\u05dc \u05dc-interpretability \u05de\u05e9\u05d5\u05e4\u05e8 \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc, \u05de\u05db\u05d9\u05d5\u05d5\u05df \u05e9\u05e0\u05d9\u05ea\u05df \u05dc\u05e0\u05ea\u05d7 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05e7\u05dc\u05d5\u05ea \u05d0\u05ea \u05d4\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d4\u05d7\u05d3-\u05de\u05e9\u05ea\u05e0\u05d9\u05d5\u05ea \u05e9\u05e0\u05dc\u05de\u05d3\u05d5.\",\n  \"\u05e8\u05e9\u05ea\u05d5\u05ea KANs \u05dc\u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd:\",\n  \"\u05db\u05e2\u05ea \u05e0\u05ea\u05d0\u05e8 \u05db\u05de\u05d4 \u05d4\u05e8\u05d7\u05d1\u05d5\u05ea \u05e9\u05dc KAN \u05dc\u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd. \u05dc\u05e0\u05d9\u05ea\u05d5\u05d7 \u05e1\u05d3\u05e8\u05d5\u05ea \u05e2\u05ea\u05d9\u05d5\u05ea, \u05e8\u05e9\u05ea\u05d5\u05ea KAN \u05d6\u05de\u05e0\u05d9\u05d5\u05ea (T-KANs) \u05de\u05e9\u05dc\u05d1\u05d5\u05ea \u05de\u05e0\u05d2\u05e0\u05d5\u05e0\u05d9 \u05d6\u05d9\u05db\u05e8\u05d5\u05df, \u05d1\u05d3\u05d5\u05de\u05d4 \u05dc-RNNs \u05d5-LSTM, \u05dc\u05d8\u05d9\u05e4\u05d5\u05dc \u05d9\u05e2\u05d9\u05dc \u05d1\u05e1\u05d3\u05e8\u05d5\u05ea \u05d0\u05dc\u05d5 \u05d5\u05d1\u05ea\u05dc\u05d5\u05d9\u05d5\u05ea \u05dc\u05d8\u05d5\u05d5\u05d7 \u05d0\u05e8\u05d5\u05da \u05e9\u05d1\u05d4\u05df, \u05d5\u05de\u05d3\u05d2\u05d9\u05de\u05d5\u05ea \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05de\u05e2\u05d5\u05dc\u05d9\u05dd \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d7\u05d9\u05d6\u05d5\u05d9 \u05e8\u05d1-\u05e9\u05dc\u05d1\u05d9(multi-step forecasting). \u05d1\u05e0\u05d5\u05e1\u05e3, \u05e9\u05d9\u05e0\u05d5\u05d9\u05d9\u05dd \u05db\u05de\u05d5 \u05de\u05e0\u05d2\u05e0\u05d5\u05e0\u05d9\u05dd \u05d7\u05d9\u05d1\u05d5\u05e8\u05d9\u05dd gated, \u05d1\u05d3\u05d5\u05de\u05d4 LSTM \u05d5-GRU, \u05de\u05d0\u05e4\u05e9\u05e8\u05d9\u05dd \u05dc-KANs \u05dc\u05d4\u05ea\u05d0\u05d9\u05dd \u05d1\u05d0\u05d5\u05e4\u05df \u05d3\u05d9\u05e0\u05de\u05d9 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d4 (\u05e1\u05e4\u05dc\u05d9\u05d9\u05df \u05d1\u05d2\u05d3\u05d5\u05dc* \u05d1\u05d4\u05ea\u05d1\u05e1\u05e1 \u05e2\u05dc \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea \u05d4\u05de\u05e9\u05d9\u05de\u05d4, \u05de\u05e9\u05e4\u05e8\u05d9\u05dd \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05de\u05d1\u05dc\u05d9 \u05dc\u05d3\u05e8\u05d5\u05e9 \u05e8\u05d2\u05d5\u05dc\u05e8\u05d9\u05d6\u05e6\u05d9\u05d4 \u05e0\u05e8\u05d7\u05d1\u05ea.\",\n  \"\u05d1\u05d3\u05d0\u05d8\u05d4 \u05d4\u05d2\u05e8\u05e4\u05d9, KANs \u05de\u05d1\u05d5\u05e1\u05e1\u05d5\u05ea \u05d2\u05e8\u05e3 (GKANs) \u05e4\u05d5\u05ea\u05d7\u05d5 \u05dc\u05e9\u05d9\u05e4\u05d5\u05e8 \u05e1\u05d9\u05d5\u05d5\u05d2 \u05e6\u05de\u05ea\u05d9\u05dd semi-supervised \u05e2\u05dc \u05d9\u05d3\u05d9 \u05e9\u05d9\u05e4\u05d5\u05e8 \u05d6\u05e8\u05d9\u05de\u05ea \u05de\u05d9\u05d3\u05e2 \u05d1\u05d9\u05df \u05e6\u05de\u05ea\u05d9\u05dd, \u05e2\u05d5\u05dc\u05d5\u05ea \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05d4\u05df \u05e2\u05dc \u05e8\u05e9\u05ea\u05d5\u05ea \u05e7\u05d5\u05e0\u05d1\u05d5\u05dc\u05d5\u05e6\u05d9\u05d4 \u05d2\u05e8\u05e4\u05d9\u05d5\u05ea \u05de\u05e1\u05d5\u05e8\u05ea\u05d9\u05d5\u05ea (GCNs). \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d5\u05ea \u05de\u05d1\u05d5\u05e1\u05e1\u05d5\u05ea KAN \u05d0\u05dc\u05d4 \u05de\u05e9\u05e4\u05e8\u05d5\u05ea \u05d0\u05ea \u05dc\u05de\u05d9\u05d3\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05e6\u05de\u05ea\u05d9\u05dd \u05d5\u05de\u05e9\u05e4\u05e8\u05d5\u05ea \u05d0\u05ea \u05d3\u05d9\u05d5\u05e7 \u05de\u05d5\u05d3\u05dc\u05d9 \u05d4\u05e8\u05d2\u05e8\u05e1\u05d9\u05d4 \u05d1\u05d2\u05e8\u05e4\u05d9\u05dd \u05d4\u05e2\u05d5\u05dc\u05d5\u05ea \u05d1\u05e8\u05e9\u05ea\u05d5\u05ea \u05d7\u05d1\u05e8\u05ea\u05d9\u05d5\u05ea \u05d5\u05db\u05d9\u05de\u05d9\u05d4 \u05de\u05d5\u05dc\u05e7\u05d5\u05dc\u05e8\u05d9\u05ea. GCNs \u05e4\u05d5\u05e2\u05dc\u05d5\u05ea \u05e2\u05dc \u05d9\u05d3\u05d9 \u05e6\u05d1\u05d9\u05e8\u05d4 \u05d5\u05e9\u05d9\u05e0\u05d5\u05d9 \u05d7\u05d5\u05d6\u05e8\u05d9\u05dd \u05e9\u05dc \u05de\u05d9\u05d3\u05e2 \u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05de\u05e9\u05db\u05d5\u05e0\u05d5\u05ea \u05de\u05e7\u05d5\u05de\u05d9\u05d5\u05ea \u05d1\u05ea\u05d5\u05da \u05d2\u05e8\u05e3, \u05d5\u05ea\u05d5\u05e4\u05e1\u05d5\u05ea \u05d1\u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d4\u05df \u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05e6\u05de\u05ea\u05d9\u05dd \u05d5\u05d4\u05df \u05d8\u05d5\u05e4\u05d5\u05dc\u05d5\u05d2\u05d9\u05d9\u05ea \u05d2\u05e8\u05e3. \u05e2\u05dd \u05d6\u05d0\u05ea, GCNs \u05de\u05e1\u05ea\u05de\u05db\u05d5\u05ea \u05e2\u05dc \u05e4\u05d9\u05dc\u05d8\u05e8\u05d9 \u05e7\u05d5\u05e0\u05d1\u05d5\u05dc\u05d5\u05e6\u05d9\u05d4 \u05e7\u05d1\u05d5\u05e2\u05d9\u05dd, \u05d4\u05de\u05d2\u05d1\u05d9\u05dc\u05d9\u05dd \u05d0\u05ea \u05d4\u05d2\u05de\u05d9\u05e9\u05d5\u05ea \u05e9\u05dc\u05d4\u05df \u05d1\u05d8\u05d9\u05e4\u05d5\u05dc \u05d1\u05d2\u05e8\u05e4\u05d9\u05dd \u05de\u05d5\u05e8\u05db\u05d1\u05d9\u05dd \u05d5\u05d4\u05d8\u05e8\u05d5\u05d2\u05e0\u05d9\u05d9\u05dd. \u05db\u05d3\u05d9 \u05dc\u05d4\u05ea\u05de\u05d5\u05d3\u05d3 \u05e2\u05dd \u05de\u05d2\u05d1\u05dc\u05d4 \u05d6\u05d5, GKAN \u05de\u05e6\u05d9\u05d2 \u05e9\u05ea\u05d9 \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d5\u05ea \u05e2\u05d9\u05e7\u05e8\u05d9\u05d5\u05ea: \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4 1, \u05d4\u05de\u05e6\u05e8\u05e4\u05ea \u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05e6\u05de\u05ea\u05d9\u05dd \u05dc\u05e4\u05e0\u05d9 \u05d9\u05d9\u05e9\u05d5\u05dd \u05e9\u05db\u05d1\u05d5\u05ea KAN, \u05de\u05d0\u05e4\u05e9\u05e8\u05ea \u05dc\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d4 \u05d4\u05e0\u05d9\u05ea\u05e0\u05d5\u05ea \u05dc\u05dc\u05de\u05d9\u05d3\u05d4 \u05dc\u05ea\u05e4\u05d5\u05e1 \u05d9\u05d7\u05e1\u05d9\u05dd \u05de\u05e7\u05d5\u05de\u05d9\u05d9\u05dd \u05de\u05d5\u05e8\u05db\u05d1\u05d9\u05dd, \u05d5\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4 2, \u05d4\u05de\u05de\u05e7\u05de\u05ea \u05e9\u05db\u05d1\u05d5\u05ea KAN \u05d1\u05d9\u05df \u05d4\u05d8\u05de\u05e2\u05d5\u05ea \u05e6\u05de\u05ea\u05d9\u05dd \u05d1\u05db\u05dc \u05e9\u05db\u05d1\u05d4 \u05dc\u05e4\u05e0\u05d9 \u05d4\u05e6\u05d1\u05d9\u05e8\u05d4, \u05de\u05d0\u05e4\u05e9\u05e8\u05ea \u05d4\u05ea\u05d0\u05de\u05d4 \u05d3\u05d9\u05e0\u05de\u05d9\u05ea \u05dc\u05e9\u05d9\u05e0\u05d5\u05d9\u05d9\u05dd \u05d1\u05de\u05d1\u05e0\u05d4 \u05d4\u05d2\u05e8\u05e3. \u05e9\u05d9\u05e4\u05d5\u05e8 \u05d6\u05d4 \u05de\u05d0\u05e4\u05e9\u05e8 \u05dc-GKANs \u05dc\u05d4\u05e1\u05ea\u05d2\u05dc \u05d1\u05d0\u05d5\u05e4\u05df \u05d3\u05d9\u05e0\u05de\u05d9 \u05dc\u05e9\u05d9\u05e0\u05d5\u05d9\u05d9\u05dd \u05d1\u05de\u05d1\u05e0\u05d4 \u05d4\u05d2\u05e8\u05e3, \u05d5\u05de\u05e1\u05e4\u05e7 \u05d2\u05d9\u05e9\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05d0\u05d3\u05e4\u05d8\u05d9\u05d1\u05d9\u05ea \u05dc\u05dc\u05de\u05d9\u05d3\u05d4 \u05de\u05d1\u05d5\u05e1\u05e1\u05ea \u05d2\u05e8\u05e3.\"\n];\n\nconst newParagraphTexts = [\n  \"\u05dc\u05e4\u05ea\u05e8\u05d5\u05df \u05de\u05e9\u05d5\u05d5\u05d0\u05d5\u05ea \u05d3\u05d9\u05e4\u05e8\u05e0\u05e6\u05d9\u05d0\u05dc\u05d9\u05d5\u05ea, KANs \u05de\u05d1\u05d5\u05e1\u05e1\u05d5\u05ea \u05e4\u05d9\u05d6\u05d9\u05e7\u05d4 (PIKANs) \u05d4\u05d5\u05ea\u05d0\u05de\u05d5 \u05dc\u05d4\u05e6\u05d9\u05e2 \u05d0\u05dc\u05d8\u05e8\u05e0\u05d8\u05d9\u05d1\u05d4 \u05e0\u05d9\u05ea\u05e0\u05ea \u05dc\u05e4\u05d9\u05e8\u05d5\u05e9(interpretability) \u05d5\u05d9\u05e2\u05d9\u05dc\u05d4 \u05dc\u05e8\u05e9\u05ea\u05d5\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd \u05de\u05d1\u05d5\u05e1\u05e1\u05d5\u05ea \u05e4\u05d9\u05d6\u05d9\u05e7\u05dc\u05d9\u05d5\u05ea (PINNs) \u05d4\u05de\u05d1\u05d5\u05e1\u05e1\u05d5\u05ea \u05e2\u05dc MLPs. \u05db\u05d0\u05df PIKANs \u05de\u05e9\u05ea\u05de\u05e9\u05d5\u05ea \u05d1\u05de\u05d1\u05e0\u05d4 \u05d0\u05d3\u05e4\u05d8\u05d9\u05d1\u05d9 \u05ea\u05dc\u05d5\u05d9-\u05d2\u05e8\u05d9\u05d3, \u05de\u05d4 \u05e9\u05d4\u05d5\u05e4\u05da \u05d0\u05d5\u05ea\u05df \u05de\u05ea\u05d0\u05d9\u05de\u05d5\u05ea \u05dc\u05d9\u05d9\u05e9\u05d5\u05de\u05d9\u05dd \u05d4\u05d3\u05d5\u05e8\u05e9\u05d9\u05dd \u05d3\u05d9\u05d5\u05e7, \u05db\u05de\u05d5 \u05d3\u05d9\u05e0\u05de\u05d9\u05e7\u05ea \u05d6\u05e8\u05d9\u05de\u05d4 \u05d5\u05de\u05db\u05e0\u05d9\u05e7\u05ea \u05e7\u05d5\u05d5\u05e0\u05d8\u05d9\u05dd, \u05e9\u05d1\u05d4\u05df \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d1\u05e1\u05d9\u05e1 \u05d3\u05d9\u05e0\u05de\u05d9\u05d5\u05ea \u05e2\u05d5\u05d6\u05e8\u05d5\u05ea \u05dc\u05ea\u05e4\u05d5\u05e1 \u05ea\u05d4\u05dc\u05d9\u05db\u05d9\u05dd \u05e4\u05d9\u05d6\u05d9\u05e7\u05dc\u05d9\u05d9\u05dd \u05de\u05d5\u05e8\u05db\u05d1\u05d9\u05dd \u05e2\u05dd \u05d3\u05d9\u05d5\u05e7 \u05d5\u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05ea \u05de\u05e9\u05d5\u05e4\u05e8\u05d9\u05dd.\",\n  \"\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d2\u05dd \u05d3\u05e0\u05d9\u05dd \u05d1\u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d4 \u05d4\u05de\u05d0\u05ea\u05d2\u05e8\u05ea \u05e9\u05dc KANs \u05d1\u05e9\u05dc \u05d4\u05d0\u05d5\u05e4\u05d9 \u05d4\u05dc\u05d0-\u05dc\u05d9\u05e0\u05d0\u05e8\u05d9 \u05e9\u05dc \u05e4\u05e8\u05de\u05d8\u05e8\u05d9 \u05d4\u05e1\u05e4\u05dc\u05d9\u05d9\u05e0\u05d9\u05dd \u05de\u05d9\u05de\u05d3\u05d9\u05d5\u05ea \u05d4\u05d2\u05d1\u05d5\u05d4\u05d4 \u05d1\u05d4 \u05e0\u05ea\u05e7\u05dc\u05d9\u05dd \u05dc\u05e2\u05d9\u05ea\u05d9\u05dd \u05e7\u05e8\u05d5\u05d1\u05d5\u05ea.\",\n  \"\u05e1\u05d9\u05db\u05d5\u05dd:\",\n  \"KANs \u05de\u05e9\u05ea\u05de\u05e9\u05d5\u05ea \u05d1-B-splines \u05dc\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05d6\u05e6\u05d9\u05d4 \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05e9\u05dc \u05de\u05e9\u05ea\u05e0\u05d4 \u05d0\u05d7\u05d3, \u05de\u05d4 \u05e9\u05d4\u05d5\u05e4\u05da \u05d0\u05d5\u05ea\u05df \u05dc\u05e0\u05d9\u05ea\u05e0\u05d5\u05ea \u05dc\u05dc\u05de\u05d9\u05d3\u05d4 \u05d5\u05de\u05d0\u05e4\u05e9\u05e8 \u05de\u05e2\u05d1\u05e8\u05d9\u05dd \u05d7\u05dc\u05e7\u05d9\u05dd \u05d1\u05d9\u05df \u05d0\u05d9\u05e0\u05d8\u05e8\u05d5\u05d5\u05dc\u05d9\u05dd \u05d4\u05e9\u05d5\u05e0\u05d9\u05dd \u05e2\u05dd \u05d4\u05ea\u05d0\u05de\u05d4 \u05de\u05e7\u05d5\u05de\u05d9\u05ea \u05de\u05e9\u05d5\u05e4\u05e8\u05ea \u05e9\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4. \u05ea\u05d4\u05dc\u05d9\u05da \u05d4\u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d4 \u05db\u05d5\u05dc\u05dc \u05d4\u05ea\u05d0\u05de\u05ea \u05e4\u05e8\u05de\u05d8\u05e8\u05d9 \u05d4\u05e1\u05e4\u05dc\u05d9\u05d9\u05e0\u05d9\u05dd, \u05db\u05de\u05d5 \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d1\u05e7\u05e8\u05d4(control point) \u05d5\u05e7\u05e9\u05e8\u05d9\u05dd, \u05db\u05d3\u05d9 \u05dc\u05de\u05d6\u05e2\u05e8 \u05e9\u05d2\u05d9\u05d0\u05d5\u05ea \u05d1\u05d9\u05df \u05e4\u05dc\u05d8 \u05d7\u05d6\u05d5\u05d9 \u05dc\u05e4\u05dc\u05d8 \u05d0\u05de\u05d9\u05ea\u05d9, \u05de\u05d0\u05e4\u05e9\u05e8 \u05dc\u05de\u05d5\u05d3\u05dc \u05dc\u05ea\u05e4\u05d5\u05e1 \u05d3\u05e4\u05d5\u05e1\u05d9 \u05d3\u05d0\u05d8\u05d4 \u05de\u05d5\u05e8\u05db\u05d1\u05d9\u05dd. \u05e2\u05dd \u05d6\u05d0\u05ea, \u05ea\u05d4\u05dc\u05d9\u05da \u05d6\u05d4 \u05de\u05e1\u05d5\u05d1\u05da \u05d1\u05e9\u05dc \u05de\u05e8\u05d7\u05d1 \u05d4\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05d4\u05dc\u05d0-\u05dc\u05d9\u05e0\u05d9\u05d0\u05e8\u05d9, \u05e7\u05dc\u05dc\u05ea \u05d4\u05de\u05de\u05d3\u05d9\u05d5\u05ea, \u05d5\u05d4\u05ea\u05e7\u05d5\u05e8\u05d4 \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05ea \u05d4\u05de\u05d5\u05d2\u05d1\u05e8\u05ea \u05d1\u05e9\u05dc \u05d4\u05d2\u05de\u05d9\u05e9\u05d5\u05ea \u05e9\u05dc \u05e1\u05e4\u05dc\u05d9\u05d9\u05e0\u05d9\u05dd \u05d4\u05e0\u05d9\u05ea\u05e0\u05d9\u05dd.\",\n  \"https://arxiv.org/abs/2411.06078\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length < existingParagraphTexts.length) {\n  throw new Error(\n    `Expected at least ${existingParagraphTexts.length} paragraphs, found ${paragraphs.items.length}`\n  );\n}\n\n// 1) Overwrite the text of the first 7 paragraphs in place.\nfor (let i = 0; i < existingParagraphTexts.length; i++) {\n  paragraphs.items[i].getRange().insertText(existingParagraphTexts[i], \"Replace\");\n}\n\nawait context.sync();\n\n// 2) Append the new paragraphs after the last (7th) paragraph, in order.\nlet cursor = paragraphs.items[existingParagraphTexts.length - 1];\nfor (const text of newParagraphTexts) {\n  cursor = cursor.insertParagraph(text, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script\n# Replaces the text of the first 7 paragraphs and appends 5 new paragraphs\n# at the end, matching the diff in the commit.\n#\n# Paragraph 1 contains a manual line break (<w:br/>) between its two runs;\n# COM represents that break as Chr(11) (vertical tab) inside Range.Text, so\n# it is rebuilt below with \"...\" + [char]11 + \"...\".\n\n$d = $word.ActiveDocument\n\n$existingParagraphTexts = @(\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 - 14.03.25\" + [char]11 + \"A Survey on Kolmogorov-Arnold Network\",\n  \"\u05de\u05d1\u05d5\u05d0:\",\n  \" \u05d6\u05d5\u05db\u05e8\u05d9\u05dd \u05d0\u05ea KANs? \u05e9\u05d6\u05d4 \u05e7\u05d9\u05e6\u05d5\u05e8 \u05e9\u05dc Kolmogorov-Arnold Networks \u05e9\u05e2\u05e9\u05d4 \u05d4\u05e8\u05d1\u05d4 \u05e8\u05e2\u05e9 \u05d1\u05d6\u05de\u05e0\u05d5 \u05d0\u05da \u05d4\u05d1\u05d0\u05d6 \u05d4\u05dc\u05da \u05d5\u05d3\u05e2\u05da \u05e2\u05dd \u05d4\u05d6\u05de\u05df. \u05de\u05ea\u05d1\u05e8\u05e8 \u05e9\u05d9\u05e6\u05d0\u05d5 \u05dc\u05d0 \u05de\u05e2\u05d8 \u05de\u05d7\u05e7\u05e8\u05d9\u05dd \u05d1\u05e0\u05d5\u05e9\u05d0 \u05d4\u05de\u05e8\u05ea\u05e7 \u05d4\u05d6\u05d4. \u05d4\u05de\u05d0\u05de\u05e8 \u05d3\u05df \u05d1\u05d4\u05e8\u05d7\u05d1\u05d5\u05ea \u05d5\u05e9\u05d9\u05e0\u05d5\u05d9\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd \u05dc\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05ea \u05d4-KAN \u05d4\u05d1\u05e1\u05d9\u05e1\u05d9\u05ea. \u05d0\u05dc\u05d4 \u05db\u05d5\u05dc\u05dc\u05d9\u05dd \u05d4\u05ea\u05d0\u05de\u05d5\u05ea \u05dc\u05e0\u05d9\u05ea\u05d5\u05d7 \u05e1\u05d3\u05e8\u05d5\u05ea \u05e2\u05ea\u05d9\u05d5\u05ea, \u05dc\u05e2\u05d9\u05d1\u05d5\u05d3 \u05d3\u05d0\u05d8\u05d4 \u05d2\u05e8\u05e4\u05d9 \u05d5\u05dc\u05e4\u05ea\u05e8\u05d5\u05df \u05de\u05e9\u05d5\u05d5\u05d0\u05d5\u05ea \u05d3\u05d9\u05e4\u05e8\u05e0\u05e6\u05d9\u05d0\u05dc\u05d9\u05d5\u05ea. \u05e9\u05d9\u05e0\u05d5\u05d9\u05d9\u05dd \u05d0\u05dc\u05d4 \u05db\u05d5\u05dc\u05dc\u05d9\u05dd \u05dc\u05e8\u05d5\u05d1 \u05e9\u05d9\u05dc\u05d5\u05d1 \u05e9\u05dc \u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05de\u05d9\u05d5\u05d7\u05d3\u05d9\u05dd \u05d0\u05d5 \u05d0\u05d9\u05dc\u05d5\u05e6\u05d9\u05dd \u05d1\u05ea\u05d5\u05da \u05d4-KAN \u05d1\u05de\u05d8\u05e8\u05d4 \u05dc\u05d4\u05ea\u05de\u05d5\u05d3\u05d3 \u05d8\u05d5\u05d1 \u05d9\u05d5\u05ea\u05e8 \u05e2\u05dd \u05d4\u05d3\u05e8\u05d9\u05e9\u05d5\u05ea \u05d4\u05e1\u05e4\u05e6\u05d9\u05e4\u05d9\u05d5\u05ea \u05e9\u05dc \u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd \u05d0\u05dc\u05d4.\",\n  \"\u05e8\u05e9\u05ea\u05d5\u05ea \u05e7\u05d5\u05dc\u05de\u05d5\u05d2\u05d5\u05e8\u05d5\u05d1-\u05d0\u05e8\u05e0\u05d5\u05dc\u05d3 \u05de\u05d9\u05d9\u05e6\u05d2\u05d5\u05ea \u05e9\u05d9\u05e0\u05d5\u05d9 \u05e4\u05e8\u05d3\u05d9\u05d2\u05de\u05d4 \u05d1\u05ea\u05db\u05e0\u05d5\u05df \u05e8\u05e9\u05ea\u05d5\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd, \u05d4\u05de\u05d1\u05d5\u05e1\u05e1\u05d5\u05ea \u05e2\u05dc \u05de\u05e2\u05d1\u05e8 \u05de\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d4 \u05e7\u05d1\u05d5\u05e2\u05d5\u05ea \u05dc\u05e7\u05e8\u05d0\u05ea \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d4\u05e0\u05d9\u05ea\u05e0\u05d5\u05ea \u05dc\u05dc\u05de\u05d9\u05d3\u05d4 \u05d4\u05e0\u05e7\u05e8\u05d0\u05d5\u05ea b-splines. \u05d4\u05d3\u05d1\u05e8 \u05e9\u05d0\u05d1 \u05d4\u05e9\u05e8\u05d0\u05d4 \u05de\u05de\u05e9\u05e4\u05d8 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05e7\u05d5\u05dc\u05de\u05d5\u05d2\u05d5\u05e8\u05d5\u05d1-\u05d0\u05e8\u05e0\u05d5\u05dc\u05d3, \u05d4\u05d8\u05d5\u05e2\u05df \u05e9\u05db\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 \u05e8\u05e6\u05d9\u05e4\u05d4 \u05e9\u05dc \u05de\u05e9\u05ea\u05e0\u05d9\u05dd \u05de\u05e8\u05d5\u05d1\u05d9\u05dd \u05e0\u05d9\u05ea\u05e0\u05ea \u05dc\u05d9\u05d9\u05e6\u05d5\u05d2 \u05db\u05d4\u05e8\u05db\u05d1\u05d4 \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05e9\u05dc \u05de\u05e9\u05ea\u05e0\u05d4 \u05d0\u05d7\u05d3. \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05e9\u05d9\u05de\u05d5\u05e9 \u05d1\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d4\u05de\u05d9\u05d5\u05e6\u05d2\u05d5\u05ea \u05e2\u05dc \u05d9\u05d3\u05d9 \u05e1\u05e4\u05dc\u05d9\u05d9\u05e0\u05d9\u05dd(\u05e9\u05d9\u05dc\u05d5\u05d1 \u05e9\u05dc \u05e4\u05d5\u05dc\u05d9\u05e0\u05d5\u05de\u05d9\u05dd \u05d1\u05d0\u05d9\u05e0\u05d8\u05e8\u05d5\u05d5\u05dc \u05e1\u05e4\u05d5\u05d9), KANs \u05de\u05e6\u05d9\u05e2\u05d5\u05ea \u05d2\u05de\u05d9\u05e9\u05d5\u05ea \u05de\u05e9\u05d5\u05e4\u05e8\u05ea \u05d5\u05e4\u05d5\u05d8\u05e0\u05e6\u05d9\u05d0\u05dc \u05dc\u05d3\u05d9\u05d5\u05e7 \u05d2\u05d1\u05d5\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05d1\u05e7\u05d9\u05e8\u05d5\u05d1 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea. \u05d3\u05d1\u05e8 \u05de\u05d5\u05d1\u05d9\u05dc \u05dc-interpretability \u05de\u05e9\u05d5\u05e4\u05e8 \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc, \u05de\u05db\u05d9\u05d5\u05d5\u05df \u05e9\u05e0\u05d9\u05ea\u05df \u05dc\u05e0\u05ea\u05d7 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05e7\u05dc\u05d5\u05ea \u05d0\u05ea \u05d4\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d4\u05d7\u05d3-\u05de\u05e9\u05ea\u05e0\u05d9\u05d5\u05ea \u05e9\u05e0\u05dc\u05de\u05d3\u05d5.\",\n  \"\u05e8\u05e9\u05ea\u05d5\u05ea KANs \u05dc\u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd:\",\n  \"\u05db\u05e2\u05ea \u05e0\u05ea\u05d0\u05e8 \u05db\u05de\u05d4 \u05d4\u05e8\u05d7\u05d1\u05d5\u05ea \u05e9\u05dc KAN \u05dc\u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd. \u05dc\u05e0\u05d9\u05ea\u05d5\u05d7 \u05e1\u05d3\u05e8\u05d5\u05ea \u05e2\u05ea\u05d9\u05d5\u05ea, \u05e8\u05e9\u05ea\u05d5\u05ea KAN \u05d6\u05de\u05e0\u05d9\u05d5\u05ea (T-KANs) \u05de\u05e9\u05dc\u05d1\u05d5\u05ea \u05de\u05e0\u05d2\u05e0\u05d5\u05e0\u05d9 \u05d6\u05d9\u05db\u05e8\u05d5\u05df, \u05d1\u05d3\u05d5\u05de\u05d4 \u05dc-RNNs \u05d5-LSTM, \u05dc\u05d8\u05d9\u05e4\u05d5\u05dc \u05d9\u05e2\u05d9\u05dc \u05d1\u05e1\u05d3\u05e8\u05d5\u05ea \u05d0\u05dc\u05d5 \u05d5\u05d1\u05ea\u05dc\u05d5\u05d9\u05d5\u05ea \u05dc\u05d8\u05d5\u05d5\u05d7 \u05d0\u05e8\u05d5\u05da \u05e9\u05d1\u05d4\u05df, \u05d5\u05de\u05d3\u05d2\u05d9\u05de\u05d5\u05ea \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05de\u05e2\u05d5\u05dc\u05d9\u05dd \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d7\u05d9\u05d6\u05d5\u05d9 \u05e8\u05d1-\u05e9\u05dc\u05d1\u05d9(multi-step forecasting). \u05d1\u05e0\u05d5\u05e1\u05e3, \u05e9\u05d9\u05e0\u05d5\u05d9\u05d9\u05dd \u05db\u05de\u05d5 \u05de\u05e0\u05d2\u05e0\u05d5\u05e0\u05d9\u05dd \u05d7\u05d9\u05d1\u05d5\u05e8\u05d9\u05dd gated, \u05d1\u05d3\u05d5\u05de\u05d4 LSTM \u05d5-GRU, \u05de\u05d0\u05e4\u05e9\u05e8\u05d9\u05dd \u05dc-KANs \u05dc\u05d4\u05ea\u05d0\u05d9\u05dd \u05d1\u05d0\u05d5\u05e4\u05df \u05d3\u05d9\u05e0\u05de\u05d9 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d4 (\u05e1\u05e4\u05dc\u05d9\u05d9\u05df \u05d1\u05d2\u05d3\u05d5\u05dc* \u05d1\u05d4\u05ea\u05d1\u05e1\u05e1 \u05e2\u05dc \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea \u05d4\u05de\u05e9\u05d9\u05de\u05d4, \u05de\u05e9\u05e4\u05e8\u05d9\u05dd \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05de\u05d1\u05dc\u05d9 \u05dc\u05d3\u05e8\u05d5\u05e9 \u05e8\u05d2\u05d5\u05dc\u05e8\u05d9\u05d6\u05e6\u05d9\u05d4 \u05e0\u05e8\u05d7\u05d1\u05ea.\",\n  \"\u05d1\u05d3\u05d0\u05d8\u05d4 \u05d4\u05d2\u05e8\u05e4\u05d9, KANs \u05de\u05d1\u05d5\u05e1\u05e1\u05d5\u05ea \u05d2\u05e8\u05e3 (GKANs) \u05e4\u05d5\u05ea\u05d7\u05d5 \u05dc\u05e9\u05d9\u05e4\u05d5\u05e8 \u05e1\u05d9\u05d5\u05d5\u05d2 \u05e6\u05de\u05ea\u05d9\u05dd semi-supervised \u05e2\u05dc \u05d9\u05d3\u05d9 \u05e9\u05d9\u05e4\u05d5\u05e8 \u05d6\u05e8\u05d9\u05de\u05ea \u05de\u05d9\u05d3\u05e2 \u05d1\u05d9\u05df \u05e6\u05de\u05ea\u05d9\u05dd, \u05e2\u05d5\u05dc\u05d5\u05ea \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05d4\u05df \u05e2\u05dc \u05e8\u05e9\u05ea\u05d5\u05ea \u05e7\u05d5\u05e0\u05d1\u05d5\u05dc\u05d5\u05e6\u05d9\u05d4 \u05d2\u05e8\u05e4\u05d9\u05d5\u05ea \u05de\u05e1\u05d5\u05e8\u05ea\u05d9\u05d5\u05ea (GCNs). \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d5\u05ea \u05de\u05d1\u05d5\u05e1\u05e1\u05d5\u05ea KAN \u05d0\u05dc\u05d4 \u05de\u05e9\u05e4\u05e8\u05d5\u05ea \u05d0\u05ea \u05dc\u05de\u05d9\u05d3\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05e6\u05de\u05ea\u05d9\u05dd \u05d5\u05de\u05e9\u05e4\u05e8\u05d5\u05ea \u05d0\u05ea \u05d3\u05d9\u05d5\u05e7 \u05de\u05d5\u05d3\u05dc\u05d9 \u05d4\u05e8\u05d2\u05e8\u05e1\u05d9\u05d4 \u05d1\u05d2\u05e8\u05e4\u05d9\u05dd \u05d4\u05e2\u05d5\u05dc\u05d5\u05ea \u05d1\u05e8\u05e9\u05ea\u05d5\u05ea \u05d7\u05d1\u05e8\u05ea\u05d9\u05d5\u05ea \u05d5\u05db\u05d9\u05de\u05d9\u05d4 \u05de\u05d5\u05dc\u05e7\u05d5\u05dc\u05e8\u05d9\u05ea. GCNs \u05e4\u05d5\u05e2\u05dc\u05d5\u05ea \u05e2\u05dc \u05d9\u05d3\u05d9 \u05e6\u05d1\u05d9\u05e8\u05d4 \u05d5\u05e9\u05d9\u05e0\u05d5\u05d9 \u05d7\u05d5\u05d6\u05e8\u05d9\u05dd \u05e9\u05dc \u05de\u05d9\u05d3\u05e2 \u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05de\u05e9\u05db\u05d5\u05e0\u05d5\u05ea \u05de\u05e7\u05d5\u05de\u05d9\u05d5\u05ea \u05d1\u05ea\u05d5\u05da \u05d2\u05e8\u05e3, \u05d5\u05ea\u05d5\u05e4\u05e1\u05d5\u05ea \u05d1\u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d4\u05df \u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05e6\u05de\u05ea\u05d9\u05dd \u05d5\u05d4\u05df \u05d8\u05d5\u05e4\u05d5\u05dc\u05d5\u05d2\u05d9\u05d9\u05ea \u05d2\u05e8\u05e3. \u05e2\u05dd \u05d6\u05d0\u05ea, GCNs \u05de\u05e1\u05ea\u05de\u05db\u05d5\u05ea \u05e2\u05dc \u05e4\u05d9\u05dc\u05d8\u05e8\u05d9 \u05e7\u05d5\u05e0\u05d1\u05d5\u05dc\u05d5\u05e6\u05d9\u05d4 \u05e7\u05d1\u05d5\u05e2\u05d9\u05dd, \u05d4\u05de\u05d2\u05d1\u05d9\u05dc\u05d9\u05dd \u05d0\u05ea \u05d4\u05d2\u05de\u05d9\u05e9\u05d5\u05ea \u05e9\u05dc\u05d4\u05df \u05d1\u05d8\u05d9\u05e4\u05d5\u05dc \u05d1\u05d2\u05e8\u05e4\u05d9\u05dd \u05de\u05d5\u05e8\u05db\u05d1\u05d9\u05dd \u05d5\u05d4\u05d8\u05e8\u05d5\u05d2\u05e0\u05d9\u05d9\u05dd. \u05db\u05d3\u05d9 \u05dc\u05d4\u05ea\u05de\u05d5\u05d3\u05d3 \u05e2\u05dd \u05de\u05d2\u05d1\u05dc\u05d4 \u05d6\u05d5, GKAN \u05de\u05e6\u05d9\u05d2 \u05e9\u05ea\u05d9 \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d5\u05ea \u05e2\u05d9\u05e7\u05e8\u05d9\u05d5\u05ea: \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4 1, \u05d4\u05de\u05e6\u05e8\u05e4\u05ea \u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05e6\u05de\u05ea\u05d9\u05dd \u05dc\u05e4\u05e0\u05d9 \u05d9\u05d9\u05e9\u05d5\u05dd \u05e9\u05db\u05d1\u05d5\u05ea KAN, \u05de\u05d0\u05e4\u05e9\u05e8\u05ea \u05dc\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d0\u05e7\u05d8\u05d9\u05d1\u05e6\u05d9\u05d4 \u05d4\u05e0\u05d9\u05ea\u05e0\u05d5\u05ea \u05dc\u05dc\u05de\u05d9\u05d3\u05d4 \u05dc\u05ea\u05e4\u05d5\u05e1 \u05d9\u05d7\u05e1\u05d9\u05dd \u05de\u05e7\u05d5\u05de\u05d9\u05d9\u05dd \u05de\u05d5\u05e8\u05db\u05d1\u05d9\u05dd, \u05d5\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4 2, \u05d4\u05de\u05de\u05e7\u05de\u05ea \u05e9\u05db\u05d1\u05d5\u05ea KAN \u05d1\u05d9\u05df \u05d4\u05d8\u05de\u05e2\u05d5\u05ea \u05e6\u05de\u05ea\u05d9\u05dd \u05d1\u05db\u05dc \u05e9\u05db\u05d1\u05d4 \u05dc\u05e4\u05e0\u05d9 \u05d4\u05e6\u05d1\u05d9\u05e8\u05d4, \u05de\u05d0\u05e4\u05e9\u05e8\u05ea \u05d4\u05ea\u05d0\u05de\u05d4 \u05d3\u05d9\u05e0\u05de\u05d9\u05ea \u05dc\u05e9\u05d9\u05e0\u05d5\u05d9\u05d9\u05dd \u05d1\u05de\u05d1\u05e0\u05d4 \u05d4\u05d2\u05e8\u05e3. \u05e9\u05d9\u05e4\u05d5\u05e8 \u05d6\u05d4 \u05de\u05d0\u05e4\u05e9\u05e8 \u05dc-GKANs \u05dc\u05d4\u05e1\u05ea\u05d2\u05dc \u05d1\u05d0\u05d5\u05e4\u05df \u05d3\u05d9\u05e0\u05de\u05d9 \u05dc\u05e9\u05d9\u05e0\u05d5\u05d9\u05d9\u05dd \u05d1\u05de\u05d1\u05e0\u05d4 \u05d4\u05d2\u05e8\u05e3, \u05d5\u05de\u05e1\u05e4\u05e7 \u05d2\u05d9\u05e9\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05d0\u05d3\u05e4\u05d8\u05d9\u05d1\u05d9\u05ea \u05dc\u05dc\u05de\u05d9\u05d3\u05d4 \u05de\u05d1\u05d5\u05e1\u05e1\u05ea \u05d2\u05e8\u05e3.\"\n)\n\n$newParagraphTexts = @(\n  \"\u05dc\u05e4\u05ea\u05e8\u05d5\u05df \u05de\u05e9\u05d5\u05d5\u05d0\u05d5\u05ea \u05d3\u05d9\u05e4\u05e8\u05e0\u05e6\u05d9\u05d0\u05dc\u05d9\u05d5\u05ea, KANs \u05de\u05d1\u05d5\u05e1\u05e1\u05d5\u05ea \u05e4\u05d9\u05d6\u05d9\u05e7\u05d4 (PIKANs) \u05d4\u05d5\u05ea\u05d0\u05de\u05d5 \u05dc\u05d4\u05e6\u05d9\u05e2 \u05d0\u05dc\u05d8\u05e8\u05e0\u05d8\u05d9\u05d1\u05d4 \u05e0\u05d9\u05ea\u05e0\u05ea \u05dc\u05e4\u05d9\u05e8\u05d5\u05e9(interpretability) \u05d5\u05d9\u05e2\u05d9\u05dc\u05d4 \u05dc\u05e8\u05e9\u05ea\u05d5\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd \u05de\u05d1\u05d5\u05e1\u05e1\u05d5\u05ea \u05e4\u05d9\u05d6\u05d9\u05e7\u05dc\u05d9\u05d5\u05ea (PINNs) \u05d4\u05de\u05d1\u05d5\u05e1\u05e1\u05d5\u05ea \u05e2\u05dc MLPs. \u05db\u05d0\u05df PIKANs \u05de\u05e9\u05ea\u05de\u05e9\u05d5\u05ea \u05d1\u05de\u05d1\u05e0\u05d4 \u05d0\u05d3\u05e4\u05d8\u05d9\u05d1\u05d9 \u05ea\u05dc\u05d5\u05d9-\u05d2\u05e8\u05d9\u05d3, \u05de\u05d4 \u05e9\u05d4\u05d5\u05e4\u05da \u05d0\u05d5\u05ea\u05df \u05de\u05ea\u05d0\u05d9\u05de\u05d5\u05ea \u05dc\u05d9\u05d9\u05e9\u05d5\u05de\u05d9\u05dd \u05d4\u05d3\u05d5\u05e8\u05e9\u05d9\u05dd \u05d3\u05d9\u05d5\u05e7, \u05db\u05de\u05d5 \u05d3\u05d9\u05e0\u05de\u05d9\u05e7\u05ea \u05d6\u05e8\u05d9\u05de\u05d4 \u05d5\u05de\u05db\u05e0\u05d9\u05e7\u05ea \u05e7\u05d5\u05d5\u05e0\u05d8\u05d9\u05dd, \u05e9\u05d1\u05d4\u05df \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05d1\u05e1\u05d9\u05e1 \u05d3\u05d9\u05e0\u05de\u05d9\u05d5\u05ea \u05e2\u05d5\u05d6\u05e8\u05d5\u05ea \u05dc\u05ea\u05e4\u05d5\u05e1 \u05ea\u05d4\u05dc\u05d9\u05db\u05d9\u05dd \u05e4\u05d9\u05d6\u05d9\u05e7\u05dc\u05d9\u05d9\u05dd \u05de\u05d5\u05e8\u05db\u05d1\u05d9\u05dd \u05e2\u05dd \u05d3\u05d9\u05d5\u05e7 \u05d5\u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05ea \u05de\u05e9\u05d5\u05e4\u05e8\u05d9\u05dd.\",\n  \"\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d2\u05dd \u05d3\u05e0\u05d9\u05dd \u05d1\u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d4 \u05d4\u05de\u05d0\u05ea\u05d2\u05e8\u05ea \u05e9\u05dc KANs \u05d1\u05e9\u05dc \u05d4\u05d0\u05d5\u05e4\u05d9 \u05d4\u05dc\u05d0-\u05dc\u05d9\u05e0\u05d0\u05e8\u05d9 \u05e9\u05dc \u05e4\u05e8\u05de\u05d8\u05e8\u05d9 \u05d4\u05e1\u05e4\u05dc\u05d9\u05d9\u05e0\u05d9\u05dd \u05de\u05d9\u05de\u05d3\u05d9\u05d5\u05ea \u05d4\u05d2\u05d1\u05d5\u05d4\u05d4 \u05d1\u05d4 \u05e0\u05ea\u05e7\u05dc\u05d9\u05dd \u05dc\u05e2\u05d9\u05ea\u05d9\u05dd \u05e7\u05e8\u05d5\u05d1\u05d5\u05ea.\",\n  \"\u05e1\u05d9\u05db\u05d5\u05dd:\",\n  \"KANs \u05de\u05e9\u05ea\u05de\u05e9\u05d5\u05ea \u05d1-B-splines \u05dc\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05d6\u05e6\u05d9\u05d4 \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05e9\u05dc \u05de\u05e9\u05ea\u05e0\u05d4 \u05d0\u05d7\u05d3, \u05de\u05d4 \u05e9\u05d4\u05d5\u05e4\u05da \u05d0\u05d5\u05ea\u05df \u05dc\u05e0\u05d9\u05ea\u05e0\u05d5\u05ea \u05dc\u05dc\u05de\u05d9\u05d3\u05d4 \u05d5\u05de\u05d0\u05e4\u05e9\u05e8 \u05de\u05e2\u05d1\u05e8\u05d9\u05dd \u05d7\u05dc\u05e7\u05d9\u05dd \u05d1\u05d9\u05df \u05d0\u05d9\u05e0\u05d8\u05e8\u05d5\u05d5\u05dc\u05d9\u05dd \u05d4\u05e9\u05d5\u05e0\u05d9\u05dd \u05e2\u05dd \u05d4\u05ea\u05d0\u05de\u05d4 \u05de\u05e7\u05d5\u05de\u05d9\u05ea \u05de\u05e9\u05d5\u05e4\u05e8\u05ea \u05e9\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4. \u05ea\u05d4\u05dc\u05d9\u05da \u05d4\u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d4 \u05db\u05d5\u05dc\u05dc \u05d4\u05ea\u05d0\u05de\u05ea \u05e4\u05e8\u05de\u05d8\u05e8\u05d9 \u05d4\u05e1\u05e4\u05dc\u05d9\u05d9\u05e0\u05d9\u05dd, \u05db\u05de\u05d5 \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d1\u05e7\u05e8\u05d4(control point) \u05d5\u05e7\u05e9\u05e8\u05d9\u05dd, \u05db\u05d3\u05d9 \u05dc\u05de\u05d6\u05e2\u05e8 \u05e9\u05d2\u05d9\u05d0\u05d5\u05ea \u05d1\u05d9\u05df \u05e4\u05dc\u05d8 \u05d7\u05d6\u05d5\u05d9 \u05dc\u05e4\u05dc\u05d8 \u05d0\u05de\u05d9\u05ea\u05d9, \u05de\u05d0\u05e4\u05e9\u05e8 \u05dc\u05de\u05d5\u05d3\u05dc \u05dc\u05ea\u05e4\u05d5\u05e1 \u05d3\u05e4\u05d5\u05e1\u05d9 \u05d3\u05d0\u05d8\u05d4 \u05de\u05d5\u05e8\u05db\u05d1\u05d9\u05dd. \u05e2\u05dd \u05d6\u05d0\u05ea, \u05ea\u05d4\u05dc\u05d9\u05da \u05d6\u05d4 \u05de\u05e1\u05d5\u05d1\u05da \u05d1\u05e9\u05dc \u05de\u05e8\u05d7\u05d1 \u05d4\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05d4\u05dc\u05d0-\u05dc\u05d9\u05e0\u05d9\u05d0\u05e8\u05d9, \u05e7\u05dc\u05dc\u05ea \u05d4\u05de\u05de\u05d3\u05d9\u05d5\u05ea, \u05d5\u05d4\u05ea\u05e7\u05d5\u05e8\u05d4 \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05ea \u05d4\u05de\u05d5\u05d2\u05d1\u05e8\u05ea \u05d1\u05e9\u05dc \u05d4\u05d2\u05de\u05d9\u05e9\u05d5\u05ea \u05e9\u05dc \u05e1\u05e4\u05dc\u05d9\u05d9\u05e0\u05d9\u05dd \u05d4\u05e0\u05d9\u05ea\u05e0\u05d9\u05dd.\",\n  \"https://arxiv.org/abs/2411.06078\"\n)\n\nif ($d.Paragraphs.Count -lt $existingParagraphTexts.Length) {\n    throw \"Expected at least $($existingParagraphTexts.Length) paragraphs, found $($d.Paragraphs.Count)\"\n}\n\n# 1) Overwrite the text of the first 7 paragraphs in place.\nfor ($i = 0; $i -lt $existingParagraphTexts.Length; $i++) {\n    $p = $d.Paragraphs.Item($i + 1)\n    $p.Range.Text = $existingParagraphTexts[$i]\n}\n\n# 2) Append the new paragraphs after the last (7th) paragraph, in order.\n$tailIndex = $existingParagraphTexts.Length\n$p = $d.Paragraphs.Item($tailIndex)\nforeach ($t in $newParagraphTexts) {\n    $p.Range.InsertParagraphAfter()\n    $tailIndex = $tailIndex + 1\n    $p = $d.Paragraphs.Item($tailIndex)\n    $p.Range.Text = $t\n}\n"}
